$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A16").Value = "LFU"
$ws.Range("B16").Value = 0.034733
$ws.Range("C16").Value = "TournamentBP"
